$d = $word.ActiveDocument

# Step 1: Replace "MEDICAL INVOICE" text with the new title text.
$d.Content.Find.Execute("MEDICAL INVOICE", $true, $false, $false, $false, $false,
                         $true, 1, $false, "SHIVALAYAS SIDDHA CLINIC INVOICE", 2)
